$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Id, Ost (Q) and Nord (R) values between row 2 and row 3
$ws.Range("A2").Value = 111659699
$ws.Range("Q2").Value = 800047
$ws.Range("R2").Value = 7239833

$ws.Range("A3").Value = 111659700
$ws.Range("Q3").Value = 799972
$ws.Range("R3").Value = 7239766

# Clear the Starttid (Z) and Sluttid (AB) columns for both rows
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
